$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G width (12.7109375 -> 11.7109375)
$ws.Columns.Item(7).ColumnWidth = 10.83

# Update numeric results for row 2 (random_forest)
$ws.Range("B2").Value = 1.9917072969183665
$ws.Range("C2").Value = 0.17949777369487804
$ws.Range("D2").Value = 1.5555997234247234
$ws.Range("E2").Value = 0.32711871149002608
$ws.Range("F2").Value = 0.57194292677681235
$ws.Range("G2").Value = 0.60328234496180966
$ws.Range("H2").Value = 0.67288128850997397
$ws.Range("I2").Value = 0.85328615764239202

# Update numeric results for row 3 (lsboost)
$ws.Range("B3").Value = 1.0819613690389691
$ws.Range("C3").Value = 0.097509135637974878
$ws.Range("D3").Value = 0.8523946312535684
$ws.Range("E3").Value = 0.096533458837287769
$ws.Range("F3").Value = 0.31069834057697793
$ws.Range("G3").Value = 0.33057002018706899
$ws.Range("H3").Value = 0.90346654116271219
$ws.Range("I3").Value = 0.95286773584691631
